# Finalized features to add every week's list
# Refresh the Trello "Board ID" / "Board Name" pairs for this week's trainers,
# and tidy the name/project separator to " - " (space-hyphen-space).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "62b9c217c880f863f19c818b"
$ws.Range("B3").Value = "Andrew Shields - Java Foundations Project"

$ws.Range("A4").Value = "62b9c218fa466f709254112e"
$ws.Range("B4").Value = "Marielle Nolasco - .NET Foundations Project"

$ws.Range("A5").Value = "62b9c2174c4c7878be9e576f"
$ws.Range("B5").Value = "Moiya Josephs - Java Foundations Project"

# Widen the default column to fit the longer entries now in the list.
$ws.StandardWidth = 206.957031

$wb.Save()
